# Updated Lino Salinas hours for L/A renewal
$wb = $excel.ActiveWorkbook

$dataform = $wb.Worksheets.Item("Dataform")
$grid     = $wb.Worksheets.Item("Grid")

# ---- Dataform: enter the new pay-period-15 (row 22) daily hours ----------
$dataform.Range("G22").Value = 8
$dataform.Range("H22").Value = 8
$dataform.Range("I22").Value = 8
$dataform.Range("J22").Value = 8
$dataform.Range("K22").Value = 8
$dataform.Range("N22").Value = 8
$dataform.Range("O22").Value = 8
$dataform.Range("P22").Value = 8
$dataform.Range("Q22").Value = 8

# ---- Dataform: add the new "Days" summary row (row 33) -------------------
$dataform.Range("D33").HorizontalAlignment = -4152
$dataform.Range("D33").Value = "Days"
$dataform.Range("E33").Formula = "=COUNT(F3:S28)"

# ---- Workbook: set the print area for Dataform ----------------------------
$dataform.PageSetup.PrintArea = "`$A`$1:`$S`$33"

# ---- Grid sheet: update the active-cell selection -------------------------
$grid.Range("I21").Select() | Out-Null

# ---- Make Dataform the active sheet / tab, with new active-cell selection -
$dataform.Activate()
$dataform.Range("H17").Select() | Out-Null
